# Commit message: "for some reason the notification is not appearing"
#
# This script reproduces the recorded edits:
#  - a new student (Dave Lee) was admitted on the "students" sheet, and the
#    placeholder row that had been sitting in "removed_students" was cleared
#    out (its row counter dropped to 0)
#  - a new student login was added on "student_pswd" (and its row counter bumped)
#  - a new teacher (Sam Davis) was added on "Teachers", along with a matching
#    teacher login on "teacher_psswd"
#  - the pending notification on "notifications" was fixed up: the topic is
#    now "COURSE ENROLLMENT" and the sender name is corrected to "Sam "
#    (previously wrongly duplicated as "Sam Davis") which is presumably why
#    the notification wasn't showing up correctly
#  - "student_pswd" became the active sheet/tab instead of "Teachers"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# students: admit the new student in row 3
# ---------------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("students")
$wsStudents.Range("A3").Value = 2
$wsStudents.Range("B3").Value = "Dave"
$wsStudents.Range("C3").Value = "Lee"
$wsStudents.Range("D3").Value = "hj"
$wsStudents.Range("E3").Value = "hj"
$wsStudents.Range("F3").Value = "hj"
$wsStudents.Range("G3").Value = "h"
$wsStudents.Range("H3").Value = "jhj"
$wsStudents.Range("J3").Value = 2
$wsStudents.Activate()
$wsStudents.Range("E7").Select()

# ---------------------------------------------------------------------
# removed_students: the placeholder row is cleared, row-count back to 0
# ---------------------------------------------------------------------
$wsRemoved = $wb.Worksheets.Item("removed_students")
$wsRemoved.Range("A2:H2").ClearContents()
$wsRemoved.Range("J3").Value = 0

# ---------------------------------------------------------------------
# student_courses: num_clubs style counter drops from 2 to 1
# ---------------------------------------------------------------------
$wsStudentCourses = $wb.Worksheets.Item("student_courses")
$wsStudentCourses.Range("L4").Value = 1
$wsStudentCourses.Activate()
$wsStudentCourses.Range("L4").Select()

# ---------------------------------------------------------------------
# student_pswd: new login added, row-count bumped, becomes active sheet
# ---------------------------------------------------------------------
$wsStudentPswd = $wb.Worksheets.Item("student_pswd")
$wsStudentPswd.Range("A3").Value = 2
$wsStudentPswd.Range("B3").Value = "student2"
$wsStudentPswd.Range("C3").Value = "Sam"
$wsStudentPswd.Range("G6").Value = 2
$wsStudentPswd.Activate()
$wsStudentPswd.Range("A5").Select()

# ---------------------------------------------------------------------
# Teachers: new teacher Sam Davis added
# ---------------------------------------------------------------------
$wsTeachers = $wb.Worksheets.Item("Teachers")
$wsTeachers.Range("A2").Value = 1
$wsTeachers.Range("B2").Value = "Sam"
$wsTeachers.Range("C2").Value = "Davis"
$wsTeachers.Range("D2").Value = "Maths Phd"
$wsTeachers.Range("E2").Value = "5 Years"
$wsTeachers.Range("H4").Value = 1
$wsTeachers.Activate()
$wsTeachers.Range("I29").Select()

# ---------------------------------------------------------------------
# teacher_psswd: new teacher login added, row-count bumped
# ---------------------------------------------------------------------
$wsTeacherPswd = $wb.Worksheets.Item("teacher_psswd")
$wsTeacherPswd.Range("A2").Value = 1
$wsTeacherPswd.Range("B2").Value = "lecturer1"
$wsTeacherPswd.Range("C2").Value = "Sam"
$wsTeacherPswd.Range("A3:C3").ClearContents()
$wsTeacherPswd.Range("F4").Value = 1
$wsTeacherPswd.Activate()
$wsTeacherPswd.Range("F4").Select()

# ---------------------------------------------------------------------
# notifications: fix topic text and the (previously wrong) sender name
# ---------------------------------------------------------------------
$wsNotifications = $wb.Worksheets.Item("notifications")
$wsNotifications.Range("B2").Value = "COURSE ENROLLMENT"
$wsNotifications.Range("D2").Value = "Sam "

# ---------------------------------------------------------------------
# Final active sheet/selection: student_pswd (matches activeTab="6")
# ---------------------------------------------------------------------
$wsStudentPswd.Activate()
$wsStudentPswd.Range("A5").Select()
